{"js": "// Commit: \"Added Restaurant component and api\"\n//\n// Appends a new \"Pipes\" section to the very end of the document body:\n//   - a Heading1 paragraph \"Pipes\"\n//   - a bulleted note (list level 0)\n//   - a bulleted sub-note (list level 1)\n//   - a bulleted note (list level 0) containing an interpolation example,\n//     with the same spell-check/grammar-check \"proofErr\" markers Word\n//     leaves around the flagged tokens (\"targetData\", \"pipeType\")\n//\n// The new paragraphs reuse the numbering definition (numId 1) already used\n// by every other bulleted paragraph in the document.\n//\n// insertParagraph()/insertText() on Word.Paragraph/Word.Range always coalesce\n// into a single <w:r> and drop <w:proofErr/> (it isn't part of the Office.js\n// object model), so to reproduce the exact run/proofErr structure we build a\n// small OOXML package -- the same \"pkg:package\" shape Range.getOoxml() /\n// insertOoxml() round-trip -- and hand it to insertOoxml(), which inserts the\n// markup largely verbatim.\n\nconst body = context.document.body;\n\n// The four new paragraphs, as plain WordprocessingML.\nconst newParagraphsXml =\n  '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t>Pipes</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:pStyle w:val=\"ListParagraph\"/>' +\n      '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    '</w:pPr>' +\n    '<w:r><w:t>They provide a way to transform values before it is displayed</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:pStyle w:val=\"ListParagraph\"/>' +\n      '<w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    '</w:pPr>' +\n    '<w:r><w:t>\\u201ctransform\\u201d you can change, modify, or reformat the values into what you want.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:pStyle w:val=\"ListParagraph\"/>' +\n      '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    '</w:pPr>' +\n    '<w:r><w:t>They are used in conjunction with interpolation so syntax \\u201c{{</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>targetData</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> | </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>pipeType</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">  }</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t>}\\u201d</w:t></w:r>' +\n  '</w:p>';\n\nconst documentXml =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n      newParagraphsXml +\n      '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>' +\n    '</w:body>' +\n  '</w:document>';\n\n// insertOoxml() expects the same \"mini package\" envelope that\n// Range.getOoxml() produces: a pkg:package with a package-relationships\n// part pointing at a word/document.xml part holding the fragment.\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>\\n' +\n  '<?mso-application progid=\"Word.Document\"?>\\n' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" ' +\n      'pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" ' +\n      'pkg:padding=\"512\">' +\n      '<pkg:xmlData>' +\n        '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n          '<Relationship Id=\"rId1\" ' +\n            'Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" ' +\n            'Target=\"word/document.xml\"/>' +\n        '</Relationships>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n      'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' + documentXml + '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\n// Insert at the very end of the body, after the last existing paragraph\n// (\"True -ng-valid\"), leaving everything before it untouched.\nbody.insertOoxml(ooxmlPackage, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Commit: \"Added Restaurant component and api\"\n#\n# Appends a new \"Pipes\" section to the very end of the document body:\n#   - a Heading1 paragraph \"Pipes\"\n#   - a bulleted note (list level 0)\n#   - a bulleted sub-note (list level 1)\n#   - a bulleted note (list level 0) containing an interpolation example,\n#     with the same spell-check/grammar-check \"proofErr\" markers Word\n#     leaves around the flagged tokens (\"targetData\", \"pipeType\")\n#\n# The new paragraphs reuse the numbering definition (numId 1) already used\n# by every other bulleted paragraph in the document, so Range.Text / plain\n# paragraph-style assignment isn't enough (it would not recreate the\n# numPr/ proofErr structure faithfully). Range.InsertXML lets us hand Word\n# the exact WordprocessingML for the new paragraphs.\n\n$d = $word.ActiveDocument\n\n$heading = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n  '<w:r><w:lastRenderedPageBreak/><w:t>Pipes</w:t></w:r>' +\n  '</w:p>'\n\n$bullet1 = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:pPr>' +\n    '<w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n  '</w:pPr>' +\n  '<w:r><w:t>They provide a way to transform values before it is displayed</w:t></w:r>' +\n  '</w:p>'\n\n$bullet2 = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:pPr>' +\n    '<w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"1\"/><w:numId w:val=\"1\"/></w:numPr>' +\n  '</w:pPr>' +\n  '<w:r><w:t>\u201ctransform\u201d you can change, modify, or reformat the values into what you want.</w:t></w:r>' +\n  '</w:p>'\n\n$bullet3 = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:pPr>' +\n    '<w:pStyle w:val=\"ListParagraph\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n  '</w:pPr>' +\n  '<w:r><w:t>They are used in conjunction with interpolation so syntax \u201c{{</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>targetData</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> | </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>pipeType</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">  }</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t>}\u201d</w:t></w:r>' +\n  '</w:p>'\n\n$xml = $heading + $bullet1 + $bullet2 + $bullet3\n\n# Collapse a range positioned at the very end of the document so the new\n# content is appended after the last existing paragraph (\"True -ng-valid\")\n# without disturbing it.\n$endRange = $d.Range()\n$endRange.Collapse(0)\n$endRange.InsertXML($xml)\n"}
